# Insert a new data row at row 76 (pushing the existing rows 76-164 down to
# 77-165), then populate the newly inserted row with its values.
#
# Source diff shows the row that used to be at position 76 (and every row
# after it, through the old row 164) simply shifting down by one row, with a
# brand-new record landing at row 76:
#   Fecha=44494, Volumen=115, Precio minimo=8000, Precio maximo=8500,
#   Precio promedio ponderado=8261, Origen="Región de Arica y Parinacota",
#   Precio $/Kg=138
# All the other fields of the new row match the pattern shared by every
# other "Berenjena" / "Femacal de La Calera" / "Coquimbo" row in this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(76).Insert()

$ws.Range("A76").Value = 3
$ws.Range("B76").Value = "Femacal de La Calera"
$ws.Range("C76").Value = "Coquimbo"
$ws.Range("D76").Value = 44494
$ws.Range("E76").Value = 5
$ws.Range("F76").Value = 100112001
$ws.Range("G76").Value = "Berenjena"
$ws.Range("H76").Value = "Sin especificar"
$ws.Range("I76").Value = "Primera"
$ws.Range("J76").Value = 115
$ws.Range("K76").Value = 8000
$ws.Range("L76").Value = 8500
$ws.Range("M76").Value = 8261
$ws.Range("N76").Value = "`$/caja 60 unidades"
$ws.Range("O76").Value = "Región de Arica y Parinacota"
$ws.Range("P76").Value = 138
$ws.Range("Q76").Value = 60
$ws.Range("R76").Value = "Hortaliza"
